$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap match data (columns F:V) between row 14 and row 15.
#    Columns A:E (index, pais, torneio, temporada, data_partida)
#    stay where they are - only the match details got reordered.
# ---------------------------------------------------------------
$row14 = $ws.Range("F14:V14").Value2
$row15 = $ws.Range("F15:V15").Value2
$ws.Range("F14:V14").Value2 = $row15
$ws.Range("F15:V15").Value2 = $row14

# ---------------------------------------------------------------
# 2) Swap match data (columns F:V) between row 42 and row 43.
# ---------------------------------------------------------------
$row42 = $ws.Range("F42:V42").Value2
$row43 = $ws.Range("F43:V43").Value2
$ws.Range("F42:V42").Value2 = $row43
$ws.Range("F43:V43").Value2 = $row42

# ---------------------------------------------------------------
# 3) Append 7 new match rows (48-54) after the existing data
#    (which ended at row 47). Formatting is copied from row 47
#    so the new rows pick up the same styles (bold index column,
#    date number format, etc.) as the rest of the table.
# ---------------------------------------------------------------
$newRowsData = @(
    ,@(47, "portugal", "liga-portugal", "2023-2024", 45192.6875, "Estoril", 2, "Vizela", 2, 2.26, "17/09/2023 15:42", 2.62, "23/09/2023 16:25", 3.41, "17/09/2023 15:42", 3.49, "23/09/2023 16:27", 3.31, "17/09/2023 15:42", 2.79, "23/09/2023 16:27", "https://www.betexplorer.com/football/portugal/liga-portugal/estoril-vizela/MROqZDjc/")
    ,@(48, "portugal", "liga-portugal", "2023-2024", 45192.6875, "Moreirense", 1, "SC Farense", 0, 2.14, "19/09/2023 06:12", 2.13, "23/09/2023 16:26", 3.55, "19/09/2023 06:12", 3.55, "23/09/2023 15:52", 3.44, "19/09/2023 06:12", 3.59, "23/09/2023 16:26", "https://www.betexplorer.com/football/portugal/liga-portugal/moreirense-sc-farense/fZBWP9Dq/")
    ,@(49, "portugal", "liga-portugal", "2023-2024", 45192.79166666666, "Casa Pia", 0, "Vitoria Guimaraes", 0, 2.47, "19/09/2023 06:12", 2.47, "23/09/2023 18:51", 3.16, "19/09/2023 06:12", 3.2, "23/09/2023 18:52", 3.15, "19/09/2023 06:12", 3.22, "23/09/2023 18:52", "https://www.betexplorer.com/football/portugal/liga-portugal/casa-pia-vitoria-guimaraes/YBV1AUDd/")
    ,@(50, "portugal", "liga-portugal", "2023-2024", 45192.89583333334, "FC Porto", 2, "Gil Vicente", 1, 1.15, "17/09/2023 15:42", 1.25, "23/09/2023 21:23", 9.33, "17/09/2023 15:42", 6.97, "23/09/2023 21:28", 16.99, "17/09/2023 15:42", 10.95, "23/09/2023 21:28", "https://www.betexplorer.com/football/portugal/liga-portugal/fc-porto-gil-vicente/EuU3Uk5S/")
    ,@(51, "portugal", "liga-portugal", "2023-2024", 45193.6875, "Chaves", 2, "Estrela", 2, 2.58, "20/09/2023 03:42", 2.67, "24/09/2023 16:28", 3.31, "20/09/2023 03:42", 3.49, "24/09/2023 16:28", 2.93, "20/09/2023 03:42", 2.73, "24/09/2023 16:28", "https://www.betexplorer.com/football/portugal/liga-portugal/chaves-estrela-da-amadora/fLWcBAbj/")
    ,@(52, "portugal", "liga-portugal", "2023-2024", 45193.79166666666, "Portimonense", 1, "Benfica", 3, 9.880000000000001, "19/09/2023 06:12", 13.34, "24/09/2023 18:59", 6.25, "19/09/2023 06:12", 7.72, "24/09/2023 18:59", 1.28, "19/09/2023 06:12", 1.2, "24/09/2023 18:51", "https://www.betexplorer.com/football/portugal/liga-portugal/portimonense-benfica/0INmYX53/")
    ,@(53, "portugal", "liga-portugal", "2023-2024", 45193.89583333334, "Braga", 4, "Boavista", 1, 1.34, "18/09/2023 20:42", 1.53, "24/09/2023 21:07", 5.67, "18/09/2023 20:42", 4.82, "24/09/2023 21:28", 8.18, "18/09/2023 20:42", 5.9, "24/09/2023 21:28", "https://www.betexplorer.com/football/portugal/liga-portugal/braga-boavista/tSVeWBzG/")
)

$startRow = 48
for ($i = 0; $i -lt $newRowsData.Count; $i++) {
    $targetRow = $startRow + $i
    $srcRange = $ws.Range("A47:V47")
    $srcRange.Copy()
    $dstRange = $ws.Range("A" + $targetRow + ":V" + $targetRow)
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats

    $rowValues = $newRowsData[$i]
    $arr = New-Object 'object[,]' 1,22
    for ($c = 0; $c -lt 22; $c++) {
        $arr[0, $c] = $rowValues[$c]
    }
    $ws.Range("A" + $targetRow + ":V" + $targetRow).Value2 = $arr
}

$excel.CutCopyMode = 0
